# Lecture 9 updates: reposition/resize several connector arrows and
# textboxes in the "list indexing" diagram on slide 3, drop four of the
# old staircase connector arrows, tweak one caption's text, and append
# one new connector arrow.
#
# NOTE on precision: this COM-interop host stores Shape.Left/Top/Width/
# Height as single-precision (32-bit) floats expressed in points, then
# truncates (rather than rounds) when converting back to EMU on save.
# A naive `emu / 12700` can therefore land 1 EMU short of the desired
# value. EmuToPt() nudges the point value upward by the smallest amount
# needed so the truncated round-trip reproduces the exact target EMU.
function EmuToPt {
    param([double]$Emu)

    $base = $Emu / 12700.0
    $hi = $base
    $step = 0.01
    while ([math]::Floor([double]([single]$hi) * 12700.0) -lt $Emu) {
        $hi = $hi + $step
    }
    $lo = $base
    for ($i = 0; $i -lt 60; $i++) {
        $mid = ($lo + $hi) / 2.0
        $got = [math]::Floor([double]([single]$mid) * 12700.0)
        if ($got -ge $Emu) {
            $hi = $mid
        } else {
            $lo = $mid
        }
    }
    return $hi
}

function SetBounds {
    param($Shape, $X, $Y, $Cx, $Cy)

    $Shape.Left = EmuToPt $X
    $Shape.Top = EmuToPt $Y
    $Shape.Width = EmuToPt $Cx
    $Shape.Height = EmuToPt $Cy
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Straight Arrow Connector 22: slide to the right-hand diagram, taller arrow.
SetBounds $s.Shapes.Item("Straight Arrow Connector 22") 6220222 3835873 809843 457200

# Straight Arrow Connector 23: same move, now flipped horizontally.
$conn23 = $s.Shapes.Item("Straight Arrow Connector 23")
SetBounds $conn23 5571613 3835873 611850 457200
$conn23.HorizontalFlip = -1

# The old "staircase" of 4 extra arrows (24-27) is no longer needed.
$s.Shapes.Item("Straight Arrow Connector 24").Delete()
$s.Shapes.Item("Straight Arrow Connector 25").Delete()
$s.Shapes.Item("Straight Arrow Connector 26").Delete()
$s.Shapes.Item("Straight Arrow Connector 27").Delete()

# Reflow the "[ elements ]" / "element , element" / "b   c" captions under
# the second (right-hand) diagram.
SetBounds $s.Shapes.Item("TextBox 58") 5422464 4213455 2681312 369332
SetBounds $s.Shapes.Item("TextBox 59") 5235442 5071749 2654986 369332

$tb60 = $s.Shapes.Item("TextBox 60")
SetBounds $tb60 5607259 5824917 2039834 369332
$tb60.TextFrame.TextRange.Text = "b                      c"

# Reposition the remaining connector arrows to match the new layout.
SetBounds $s.Shapes.Item("Straight Arrow Connector 61") 7010256 5411116 0 457200
SetBounds $s.Shapes.Item("Straight Arrow Connector 64") 6279769 4557629 750296 514120
SetBounds $s.Shapes.Item("Straight Arrow Connector 66") 5571613 4557629 707090 514120
SetBounds $s.Shapes.Item("Straight Arrow Connector 67") 6278238 4584765 0 457200

$conn68 = $s.Shapes.Item("Straight Arrow Connector 68")
SetBounds $conn68 5748170 5411116 0 457200

# New arrow, cloned from connector 68 so it inherits the same line/arrow
# style, then moved into place and renamed to continue the numbering.
$newConnRange = $conn68.Duplicate()
$newConn = $newConnRange.Item(1)
$newConn.Name = "Straight Arrow Connector 69"
SetBounds $newConn 6220222 3835873 0 457200
